$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.253.30'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '2.098.99'
$ws.Range('E3').Value = '  +2.94%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'229.87"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').Value = "'61.14"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.63%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = "'0.381"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('D10').Value = "'0.0842"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').Value = '2.408.41'
$ws.Range('E12').Value = '  +2.85%  '
$ws.Range('D13').Value = "'22.49"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.90%  '
$ws.Range('D14').Value = "'14.71"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.25%  '
$ws.Range('D15').Value = "'5.51"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.16%  '
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('D17').Value = '2.092.14'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').Value = '38.145.75'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = "'6.02"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('D20').Value = "'70.32"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('D22').Value = "'224.10"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('E25').Value = '  +3.33%  '
$ws.Range('D26').Value = "'170.19"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('D27').Value = "'9.45"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('D29').Value = "'19.08"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('E30').Value = '  +4.88%  '
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  +9.81%  '
$ws.Range('E33').Value = '  +3.38%  '
$ws.Range('D34').Value = "'4.44"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').Value = "'0.0607"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').Value = "'6.54"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('E37').Value = '  +5.13%  '
$ws.Range('D38').Value = "'3.56"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.69%  '
$ws.Range('D39').Value = "'0.999"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').Value = "'18.05"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('D41').Value = '1.548.39'
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('D42').Value = "'100.12"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.77%  '
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').Value = "'0.0907"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').Value = "'1.12"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.84%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = "'4.11"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('E48').Value = '  +1.69%  '
$ws.Range('E49').Value = '  +2.02%  '
$ws.Range('D51').Value = '2.296.09'
$ws.Range('E51').Value = '  +2.95%  '
